$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to Text format so numeric-looking
# strings (e.g. "0.488", "210.34") are stored as text, matching the
# workbook's existing inline-string convention for column D.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '28.684.77'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '1.565.00'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').Value = '210.34'
$ws.Range('D6').Value = '0.488'
$ws.Range('E6').Value = '  -0.51%  '
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.48%  '
$ws.Range('D8').Value = '25.18'
$ws.Range('E8').Value = '  +5.79%  '
$ws.Range('D9').Value = '0.245'
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('D11').Value = '0.0895'
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = '1.788.19'
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('D13').Value = '1.569.64'
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('D14').Value = '28.676.94'
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('D15').Value = '0.517'
$ws.Range('E15').Value = '  +0.81%  '
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').Value = '228.81'
$ws.Range('E18').Value = '  +0.71%  '
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('D20').Value = '0.0₃0679'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('D22').Value = '3.92'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').Value = '9.04'
$ws.Range('E23').Value = '  +1.08%  '
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('D25').Value = '151.45'
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('D26').Value = '14.77'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('E27').Value = '  +0.64%  '
$ws.Range('E28').Value = '  -0.39%  '
$ws.Range('E29').Value = '  -1.49%  '
$ws.Range('D30').Value = '0.0461'
$ws.Range('E30').Value = '  -3.91%  '
$ws.Range('E31').Value = '  -2.24%  '
$ws.Range('D32').Value = '3.18'
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('D33').Value = '1.394.16'
$ws.Range('E33').Value = '  +1.20%  '
$ws.Range('E34').Value = '  -2.96%  '
$ws.Range('E35').Value = '  -4.35%  '
$ws.Range('E36').Value = '  -1.09%  '
$ws.Range('E37').Value = '  +1.43%  '
$ws.Range('E38').Value = '  -2.24%  '
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('E40').Value = '  +1.50%  '
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('E42').Value = '  -0.39%  '
$ws.Range('D43').Value = '0.772'
$ws.Range('E43').Value = '  -0.88%  '
$ws.Range('E44').Value = '  -3.54%  '
$ws.Range('D45').Value = '64.07'
$ws.Range('E45').Value = '  +3.07%  '
$ws.Range('E46').Value = '  -1.63%  '
$ws.Range('D47').Value = '1.700.45'
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('D48').Value = '0.871'
$ws.Range('E48').Value = '  -5.02%  '
$ws.Range('D49').Value = '85.12'
$ws.Range('D50').Value = '43.33'
$ws.Range('E50').Value = '  +7.60%  '
$ws.Range('E51').Value = '  -0.70%  '
